$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) column titles
$ws.Cells.Item(1,1).Value = "mx_state"
$ws.Cells.Item(1,2).Value = "mx_municipality"
$ws.Cells.Item(1,3).Value = "n_matriculas"
$ws.Cells.Item(1,4).Value = "pct_matriculas"

# 2. Title-case Spanish connector words ("de","del","el","la","los","las","y")
#    inside the state (col A) and municipality (col B) name strings, for all
#    data rows (2 through 1146). The first word of each string is left as-is.
$connectors = @("de","del","el","la","los","las","y")
$lastDataRow = 1146

for ($r = 2; $r -le $lastDataRow; $r++) {
    for ($col = 1; $col -le 2; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val -ne "") {
            $words = $val.Split(" ")
            if ($words.Length -gt 1) {
                for ($i = 1; $i -lt $words.Length; $i++) {
                    $lw = $words[$i].ToLower()
                    if ($connectors -contains $lw) {
                        $cap = $lw.Substring(0,1).ToUpper() + $lw.Substring(1)
                        $words[$i] = $cap
                    }
                }
                $newval = [string]::Join(" ", $words)
                $cell.Value = $newval
            }
        }
    }
}

# 3. Fix the final grand-total label: "TOTAL" -> "Total"
$ws.Cells.Item($lastDataRow, 1).Value = "Total"

# 4. Remove the trailing metadata/footnote rows (old rows 1148-1152) that are
#    no longer part of the cleaned dataset.
$ws.Rows("1148:1152").Delete()
